$d = $word.ActiveDocument

# Locate the anchor paragraph: the one whose trimmed text is exactly nine dashes,
# immediately following the log-transformation conclusion paragraph.
$count = $d.Paragraphs.Count
$anchorIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text.Trim() -eq "---------") {
        $anchorIdx = $i
    }
}
if ($anchorIdx -eq -1) {
    throw "Could not find anchor paragraph (---------)"
}

$curIdx = $anchorIdx

# --- new paragraph #0 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
# left blank: InsertParagraphAfter already produced the empty paragraph

# --- new paragraph #1 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
$newPara = $d.Paragraphs.Item($curIdx)
$newPara.Range.Text = "A correlation matrix is a table showing correlation coefficients between sets of variables. Each random variable (Xi) in the table is correlated with each of the other values in the table (Xj). This allows you to see which pairs have the highest correlation."

# --- new paragraph #2 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
# left blank: InsertParagraphAfter already produced the empty paragraph

# --- new paragraph #3 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
$newPara = $d.Paragraphs.Item($curIdx)
$newPara.Range.Text = "This matrix provides pairwise correlation values for the selected variables. The values range from -1 to 1, with -1 indicating a perfect negative correlation, 1 indicating a perfect positive correlation, and 0 indicating no correlation."

# --- new paragraph #4 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
# left blank: InsertParagraphAfter already produced the empty paragraph

# --- new paragraph #5 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
$newPara = $d.Paragraphs.Item($curIdx)
$newPara.Range.Text = "Highly Correlated Variables:"

# --- new paragraph #6 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
# left blank: InsertParagraphAfter already produced the empty paragraph

# --- new paragraph #7 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
$newPara = $d.Paragraphs.Item($curIdx)
$newPara.Range.Text = "power.GUESS..GW, annualised.consumption.GUESS..TWh, Estimated..MtCO2e, Coal.only..MtCO2e, and Hash.rate.MH.s are highly positively correlated with each other. This means as one of these variables increases, the others tend to increase as well."

# --- new paragraph #8 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
$newPara = $d.Paragraphs.Item($curIdx)
$newPara.Range.Text = "This could indicate that as the estimated power used in bitcoin mining increases, the estimated carbon emissions and the hash rate (a measure of mining computational power) also tend to increase. This suggests a direct link between energy consumption, carbon footprint, and the computational power of the bitcoin network."

# --- new paragraph #9 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
# left blank: InsertParagraphAfter already produced the empty paragraph

# --- new paragraph #10 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
$newPara = $d.Paragraphs.Item($curIdx)
$newPara.Range.Text = "Emission Intensities:"

# --- new paragraph #11 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
# left blank: InsertParagraphAfter already produced the empty paragraph

# --- new paragraph #12 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
$newPara = $d.Paragraphs.Item($curIdx)
$newPara.Range.Text = "Emission.intensity..gCO2e.kWh does not have strong correlations with most of the other variables. This might suggest that the intensity of emissions (emissions per unit of energy) remains relatively constant regardless of the other fluctuations in the network."

# --- new paragraph #13 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
# left blank: InsertParagraphAfter already produced the empty paragraph

# --- new paragraph #14 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
$newPara = $d.Paragraphs.Item($curIdx)
$newPara.Range.Text = "Efficiency:"

# --- new paragraph #15 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
# left blank: InsertParagraphAfter already produced the empty paragraph

# --- new paragraph #16 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
$newPara = $d.Paragraphs.Item($curIdx)
$newPara.Range.Text = "The Estimated.efficiency..J.Th variable is weakly correlated with most of the other variables. This suggests that the efficiency of mining hardware (in terms of energy consumed per transaction) might not be a dominant factor influencing the overall energy consumption and emissions of the bitcoin network. However, it's notable that efficiency hasn't drastically improved or worsened significantly over time."

# --- new paragraph #17 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
# left blank: InsertParagraphAfter already produced the empty paragraph

# --- new paragraph #18 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
$newPara = $d.Paragraphs.Item($curIdx)
$newPara.Range.Text = "Hydro vs Coal Emissions:"

# --- new paragraph #19 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
# left blank: InsertParagraphAfter already produced the empty paragraph

# --- new paragraph #20 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
$newPara = $d.Paragraphs.Item($curIdx)
$newPara.Range.Text = "Hydro.only..MtCO2e has weak correlations with the other variables, indicating that emissions from hydroelectric sources don't play a significant role in the overall carbon footprint of bitcoin mining."

# --- new paragraph #21 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
$newPara = $d.Paragraphs.Item($curIdx)
$newPara.Range.Text = "On the other hand, Coal.only..MtCO2e is highly correlated with overall estimated emissions (Estimated..MtCO2e). This suggests that coal-based power sources might be a significant contributor to bitcoin's carbon footprint."

# --- new paragraph #22 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
# left blank: InsertParagraphAfter already produced the empty paragraph

# --- new paragraph #23 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
# left blank: InsertParagraphAfter already produced the empty paragraph

# --- new paragraph #24 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
$newPara = $d.Paragraphs.Item($curIdx)
$newPara.Range.Text = "## Spatial Analysis"

# --- new paragraph #25 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
# left blank: InsertParagraphAfter already produced the empty paragraph

# --- new paragraph #26 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
# left blank: InsertParagraphAfter already produced the empty paragraph

# --- new paragraph #27 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
$newPara = $d.Paragraphs.Item($curIdx)
$newPara.Range.Text = "This plot visualizes the monthly absolute hashrate for different countries using a choropleth map, a type of thematic map where the areas shaded based on the value of a variable. from looking at the map we can say that hashrate isn't evenly distributed globally but is concentrated in specific regions or countries due to various factors like technology infrastructure, regulations or energy costs.The colour gradiant ranges from light blue (reprenting lower values) to darkblue (representing higher values) which is used to represent the monthly absolute hashrate value for each country. Any country without data(Na values) is colored in grey. The countries which are shaded dark blue, it would suggest that region as a whole is a significant contributor to hashrate. The countries with higher hashrates might be inferred to have a more significant investment or infrastructure related to the domain."

# --- new paragraph #28 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
$newPara = $d.Paragraphs.Item($curIdx)
$pos = $newPara.Range.Start
$mainText = "As of Jan 2022, the maximum contribution in monthly hash rate was of United States followed by China. "
$lb = [string][char]11
$ip1 = $d.Range($pos, $pos)
$ip1.InsertAfter($lb + $mainText)
$ip2 = $d.Range($pos, $pos)
$ip2.InsertBreak(6)

# --- new paragraph #29 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
$newPara = $d.Paragraphs.Item($curIdx)
$newPara.Range.Text = "Except for US and China, most countries have had a monthly hash rate contribution of 15% or less. "

# --- new paragraph #30 ---
$curPara = $d.Paragraphs.Item($curIdx)
$curPara.Range.InsertParagraphAfter()
$curIdx = $curIdx + 1
$newPara = $d.Paragraphs.Item($curIdx)
$newPara.Range.Text = "This answers that mining activity varies with regional variations and with the countries listed, being the top contributors. "
